$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that gets bumped by one day
# for every data row (rows 2 through 465) as part of the automatic update.
$lastRow = 465
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = 46076
}
